$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "71+8=79"
$t.Cell(1, 2).Range.Text = "61-46=15"
$t.Cell(1, 3).Range.Text = "27+13=40"
$t.Cell(1, 4).Range.Text = "84-31=53"
$t.Cell(1, 5).Range.Text = "62-30=32"
$t.Cell(2, 1).Range.Text = "24+55=79"
$t.Cell(2, 2).Range.Text = "11+76=87"
$t.Cell(2, 3).Range.Text = "80-59=21"
$t.Cell(2, 4).Range.Text = "70-48=22"
$t.Cell(2, 5).Range.Text = "6+41=47"
$t.Cell(3, 1).Range.Text = "39+16=55"
$t.Cell(3, 2).Range.Text = "94-26=68"
$t.Cell(3, 3).Range.Text = "2+50=52"
$t.Cell(3, 4).Range.Text = "77+17=94"
$t.Cell(3, 5).Range.Text = "92-24=68"
$t.Cell(4, 1).Range.Text = "81-37=44"
$t.Cell(4, 2).Range.Text = "25+26=51"
$t.Cell(4, 3).Range.Text = "79+19=98"
$t.Cell(4, 4).Range.Text = "94-59=35"
$t.Cell(4, 5).Range.Text = "29+41=70"
$t.Cell(5, 1).Range.Text = "49-8=41"
$t.Cell(5, 2).Range.Text = "13+30=43"
$t.Cell(5, 3).Range.Text = "73+21=94"
$t.Cell(5, 4).Range.Text = "24+49=73"
$t.Cell(5, 5).Range.Text = "16+34=50"
$t.Cell(6, 1).Range.Text = "36+11=47"
$t.Cell(6, 2).Range.Text = "42+50=92"
$t.Cell(6, 3).Range.Text = "62-51=11"
$t.Cell(6, 4).Range.Text = "78-49=29"
$t.Cell(6, 5).Range.Text = "52-3=49"
$t.Cell(7, 1).Range.Text = "8+58=66"
$t.Cell(7, 2).Range.Text = "52-45=7"
$t.Cell(7, 3).Range.Text = "27-7=20"
$t.Cell(7, 4).Range.Text = "86-58=28"
$t.Cell(7, 5).Range.Text = "65-25=40"
$t.Cell(8, 1).Range.Text = "30+37=67"
$t.Cell(8, 2).Range.Text = "26+41=67"
$t.Cell(8, 3).Range.Text = "58-30=28"
$t.Cell(8, 4).Range.Text = "37+27=64"
$t.Cell(8, 5).Range.Text = "98-46=52"
$t.Cell(9, 1).Range.Text = "12-8=4"
$t.Cell(9, 2).Range.Text = "32+2=34"
$t.Cell(9, 3).Range.Text = "51+24=75"
$t.Cell(9, 4).Range.Text = "35+34=69"
$t.Cell(9, 5).Range.Text = "86-67=19"
$t.Cell(10, 1).Range.Text = "68-2=66"
$t.Cell(10, 2).Range.Text = "33-15=18"
$t.Cell(10, 3).Range.Text = "76+10=86"
$t.Cell(10, 4).Range.Text = "80-5=75"
$t.Cell(10, 5).Range.Text = "75+6=81"
$t.Cell(11, 1).Range.Text = "86+1=87"
$t.Cell(11, 2).Range.Text = "94-22=72"
$t.Cell(11, 3).Range.Text = "75-42=33"
$t.Cell(11, 4).Range.Text = "76-2=74"
$t.Cell(11, 5).Range.Text = "75-14=61"
$t.Cell(12, 1).Range.Text = "55+19=74"
$t.Cell(12, 2).Range.Text = "65+19=84"
$t.Cell(12, 3).Range.Text = "95-63=32"
$t.Cell(12, 4).Range.Text = "70-44=26"
$t.Cell(12, 5).Range.Text = "33+64=97"
$t.Cell(13, 1).Range.Text = "37+40=77"
$t.Cell(13, 2).Range.Text = "99-18=81"
$t.Cell(13, 3).Range.Text = "48-39=9"
$t.Cell(13, 4).Range.Text = "21+43=64"
$t.Cell(13, 5).Range.Text = "92-67=25"
$t.Cell(14, 1).Range.Text = "10+28=38"
$t.Cell(14, 2).Range.Text = "66-64=2"
$t.Cell(14, 3).Range.Text = "62-20=42"
$t.Cell(14, 4).Range.Text = "20+44=64"
$t.Cell(14, 5).Range.Text = "86-38=48"
$t.Cell(15, 1).Range.Text = "72-6=66"
$t.Cell(15, 2).Range.Text = "97-62=35"
$t.Cell(15, 3).Range.Text = "31+60=91"
$t.Cell(15, 4).Range.Text = "75-37=38"
$t.Cell(15, 5).Range.Text = "74-42=32"
$t.Cell(16, 1).Range.Text = "55-11=44"
$t.Cell(16, 2).Range.Text = "79+18=97"
$t.Cell(16, 3).Range.Text = "55+41=96"
$t.Cell(16, 4).Range.Text = "25+59=84"
$t.Cell(16, 5).Range.Text = "24+75=99"
$t.Cell(17, 1).Range.Text = "32+20=52"
$t.Cell(17, 2).Range.Text = "49-22=27"
$t.Cell(17, 3).Range.Text = "86-1=85"
$t.Cell(17, 4).Range.Text = "97-19=78"
$t.Cell(17, 5).Range.Text = "23+13=36"
$t.Cell(18, 1).Range.Text = "13-9=4"
$t.Cell(18, 2).Range.Text = "19+21=40"
$t.Cell(18, 3).Range.Text = "51-34=17"
$t.Cell(18, 4).Range.Text = "73-53=20"
$t.Cell(18, 5).Range.Text = "35+2=37"
$t.Cell(19, 1).Range.Text = "0+28=28"
$t.Cell(19, 2).Range.Text = "50+11=61"
$t.Cell(19, 3).Range.Text = "58+0=58"
$t.Cell(19, 4).Range.Text = "10+86=96"
$t.Cell(19, 5).Range.Text = "11-9=2"
$t.Cell(20, 1).Range.Text = "44+27=71"
$t.Cell(20, 2).Range.Text = "61-52=9"
$t.Cell(20, 3).Range.Text = "3+76=79"
$t.Cell(20, 4).Range.Text = "69-30=39"
$t.Cell(20, 5).Range.Text = "25+70=95"
